$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "test_summary"

# 2. Collapse the gap between Day 1 and Day 2 blocks (delete one blank row)
$ws.Rows(6).Delete()

# After the delete:
#   Day1 header/data: rows 2-5 (unchanged)
#   Day2 header/data: rows 7-10
#   Day3 header/data: rows 12-15

# 3. Update Day 3 values (content change, not just a shift)
$ws.Range("C13").Value = 121
$ws.Range("C14").Value = 79
$ws.Range("C15").Value = 167

# 4. Build the Day 4 / Day 5 / Day 6 blocks by replicating the formatting
#    of the existing Day 2 block (rows 7:10) via a formats-only paste, then
#    filling in the text/values. Text is filled in day order (4,5,6) so the
#    new shared strings land in that order; the header cells are merged
#    afterwards (day 6 first) to match the source authoring order.

function Fill-DayBlock($headerRow, $headerText, $written, $execution, $review) {
    $ws.Range("B7:C10").Copy()
    $ws.Range("B$headerRow").PasteSpecial(-4122)

    $dataRow1 = $headerRow + 1
    $dataRow2 = $headerRow + 2
    $dataRow3 = $headerRow + 3

    $ws.Range("B$headerRow").Value = $headerText
    $ws.Range("B$dataRow1").Value = "Total  testcase Written"
    $ws.Range("C$dataRow1").Value = $written
    $ws.Range("B$dataRow2").Value = "Total Execution"
    $ws.Range("C$dataRow2").Value = $execution
    $ws.Range("B$dataRow3").Value = "Total Review"
    $ws.Range("C$dataRow3").Value = $review

    $ws.Range("B${headerRow}:C${dataRow3}").RowHeight = 18
}

Fill-DayBlock 17 "Spint( 35) - Day 4 - Test Case Summary" 143 101 189
Fill-DayBlock 23 "Spint( 35) - Day 5 - Test Case Summary" 192 109 243
Fill-DayBlock 28 "Spint( 35) - Day 6 - Test Case Summary" 223 109 274

# Merge the new header rows -- day 6 first, then day 4, then day 5, matching
# the order the merges appear in the saved workbook.
$ws.Range("B28:C28").Merge()
$ws.Range("B17:C17").Merge()
$ws.Range("B23:C23").Merge()

# 5. Column C width
$ws.Columns("C").ColumnWidth = 18.17

# 6. Selection
$ws.Range("C29").Select()
